# Update "paises.xlsx" COVID-19 stats & the "last updated" timestamp.
# Sheet "Pais" holds a table sorted descending by column B ("Casos totales").
# Two pairs of neighbouring rows swap rank because of the refreshed figures:
#   - row 15 (was Canada) / row 16 (was Belgica)  -> Belgica now outranks Canada
#   - row 129 (was Islas Feroe) / row 130 (was Congo) -> Congo now outranks Islas Feroe
# For those two pairs we therefore also update column A (country name) in
# addition to the numeric columns; for every other affected row only the
# numeric columns change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Footer timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 01:52"

# --- Row 4: Estados Unidos --------------------------------------------
$ws.Range("B4").Value = 923612
$ws.Range("C4").Value = 37170
$ws.Range("D4").Value = 93404
$ws.Range("E4").Value = 778116
$ws.Range("G4").Value = 1858
$ws.Range("H4").Value = 52092

# --- Row 14: Brasil ----------------------------------------------------
$ws.Range("D14").Value = 27655
$ws.Range("E14").Value = 21670

# --- Rows 15/16: Canada & Belgica swap rank -----------------------------
# Row 15 becomes Belgica (higher "casos totales" than the refreshed Canada)
$ws.Range("A15").Value = "Belgica"
$ws.Range("B15").Value = 44293
$ws.Range("C15").Value = 1496
$ws.Range("D15").Value = 10122
$ws.Range("E15").Value = 27492
$ws.Range("F15").Value = 970
$ws.Range("G15").Value = 189
$ws.Range("H15").Value = 6679

# Row 16 becomes Canada
$ws.Range("A16").Value = "Canada"
$ws.Range("B16").Value = 43888
$ws.Range("C16").Value = 1778
$ws.Range("D16").Value = 15469
$ws.Range("E16").Value = 26117
$ws.Range("F16").Value = 557
$ws.Range("G16").Value = 155
$ws.Range("H16").Value = 2302

# --- Row 44: Chequia -----------------------------------------------------
$ws.Range("B44").Value = 7273
$ws.Range("C44").Value = 86
$ws.Range("E44").Value = 4688

# --- Row 50: Colombia ------------------------------------------------------
$ws.Range("F50").Value = 117

# --- Row 87: Hong Kong -----------------------------------------------------
$ws.Range("D87").Value = 725
$ws.Range("E87").Value = 307
$ws.Range("F87").Value = 7

# --- Row 101: Burkina Faso ---------------------------------------------
$ws.Range("B101").Value = 629
$ws.Range("C101").Value = 13
$ws.Range("D101").Value = 425
$ws.Range("E101").Value = 163

# --- Row 124: Vietnam ----------------------------------------------------
$ws.Range("D124").Value = 220
$ws.Range("E124").Value = 50

# --- Rows 129/130: Islas Feroe & Congo swap rank --------------------------
# Row 129 becomes Congo (higher refreshed "casos totales")
$ws.Range("A129").Value = "Congo"
$ws.Range("B129").Value = 200
$ws.Range("C129").Value = 14
$ws.Range("D129").Value = 19
$ws.Range("E129").Value = 175
$ws.Range("H129").Value = 6

# Row 130 becomes Islas Feroe
$ws.Range("A130").Value = "Islas Feroe"
$ws.Range("B130").Value = 187
$ws.Range("D130").Value = 178
$ws.Range("E130").Value = 9
$ws.Range("H130").Value = 0

# --- Row 146: Aruba ------------------------------------------------------
$ws.Range("D146").Value = 69
$ws.Range("E146").Value = 29

# --- Row 160: Islas Caimanes ---------------------------------------------
$ws.Range("B160").Value = 70
$ws.Range("C160").Value = 4
$ws.Range("D160").Value = 8
$ws.Range("E160").Value = 61

# --- Row 163: Polinesia Francesa ------------------------------------------
$ws.Range("D163").Value = 41
$ws.Range("E163").Value = 16

# --- Row 172: San Martin (Parte Francesa) ---------------------------------
$ws.Range("D172").Value = 24
$ws.Range("E172").Value = 11
$ws.Range("G172").Value = 1
$ws.Range("H172").Value = 3

# --- Row 175: Malaui -------------------------------------------------------
$ws.Range("D175").Value = 4
$ws.Range("E175").Value = 26
